{"js": "// Insert \" e aggiornato con delibera n\u00b0 177/2025, Verb. 521\" right before the\n// trailing \";\" of the sentence that cites delibera n\u00b0 137/2017, in the\n// \"VISTO il Codice di comportamento ...\" paragraph.\n\nconst body = context.document.body;\n\n// Anchor on the unique run of text that ends with \"n\u00b0 137/2017;\" so we can\n// locate the exact insertion point (just before the final semicolon).\nconst searchResults = body.search(\"Amministrazione n\u00b0 137/2017\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find the target text \"Amministrazione n\u00b0 137/2017\" in the document.');\n}\n\nconst target = searchResults.items[0];\nconst insertionPoint = target.getRange(\"End\");\ninsertionPoint.insertText(\" e aggiornato con delibera n\u00b0 177/2025, Verb. 521\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Insert \" e aggiornato con delibera n\u00b0 177/2025, Verb. 521\" right before the\n# trailing \";\" of the sentence that cites delibera n\u00b0 137/2017, in the\n# \"VISTO il Codice di comportamento ...\" paragraph.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \"Amministrazione n\u00b0 137/2017;\"\n$find.Replacement.Text = \"Amministrazione n\u00b0 137/2017 e aggiornato con delibera n\u00b0 177/2025, Verb. 521;\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n$find.Wrap = 1\n\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)\n"}
